$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 holds a numeric-looking identifier that must stay text; use a scratch
# cell formatted as Text, copy/paste-special the value in, then fully clear
# the scratch cell (formatting included) so the sheet dimensions/styles are
# left exactly as before.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "@"
$scratch.Value = "123"
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues: copy the text value only, not the scratch formatting
$scratch.Clear()

$ws.Range("C2").Value = [double]"0.718586444520514"
$ws.Range("D2").Value = [double]"29.6746496300906"
$ws.Range("E2").Value = [double]"26.06352488665569"
$ws.Range("F2").Value = [double]"5.619822827046354"
$ws.Range("G2").Value = [double]"1.870709232782207"
$ws.Range("H2").Value = [double]"0.3266931037974307"
$ws.Range("I2").Value = [double]"35.64843764405722"
$ws.Range("J2").Value = [double]"0.06136646676648708"
$ws.Range("K2").Value = [double]"0.01603810455514962"
$ws.Range("M2").Value = [double]"2.89664482583007e-05"
$ws.Range("N2").Value = [double]"9.023684867224475e-07"
$ws.Range("O2").Value = [double]"2.992868385413704e-09"
$ws.Range("P2").Value = [double]"1.107957643139213e-14"
$ws.Range("Q2").Value = [double]"3.456977190091634e-15"
$ws.Range("R2").Value = [double]"2.335856974189808e-17"
$ws.Range("S2").Value = [double]"1.369152556427783e-18"
$ws.Range("T2").Value = [double]"6.962386938276468e-19"
$ws.Range("U2").Value = [double]"3.642610182210246e-22"
$ws.Range("V2").Value = [double]"7.346244449418279e-20"
$ws.Range("X2").Value = [double]"2.411397523120941e-22"
$ws.Range("Y2").Value = [double]"6.576320922633137e-11"
$ws.Range("Z2").Value = [double]"5.680634416134499e-32"
$ws.Range("AB2").Value = [double]"6.083523919366755e-09"
$ws.Range("AD2").Value = [double]"3.932354624384566e-09"
$ws.Range("AE2").Value = [double]"1.218613110459088e-27"
$ws.Range("AF2").Value = [double]"1.295874703398041e-10"
$ws.Range("AG2").Value = [double]"3.498663646642182e-12"
$ws.Range("AH2").Value = [double]"0.0001213577154757903"
$ws.Range("AI2").Value = [double]"2.041997214377556e-05"
$ws.Range("AJ2").Value = [double]"1.619609257926115e-11"
$ws.Range("AL2").Value = [double]"1.756904205402146e-13"
$ws.Range("AM2").Value = [double]"7.64265530312112e-15"
$ws.Range("AN2").Value = [double]"5.18556880748865e-18"
$ws.Range("AR2").Value = [double]"0.001023705649120322"
$ws.Range("AS2").Value = [double]"0.2960200293145421"
$ws.Range("AT2").Value = [double]"0.9683346880231811"
$ws.Range("AU2").Value = [double]"0.7758658791209584"
$ws.Range("AV2").Value = [double]"0.7493582354081466"
$ws.Range("AW2").Value = [double]"0.4262177474310581"
$ws.Range("AX2").Value = [double]"96.6193798343571"
$ws.Range("AY2").Value = [double]"0.1159460796151261"
$ws.Range("AZ2").Value = [double]"0.04768917583221836"
$ws.Range("BB2").Value = [double]"0.0001554526653199358"
$ws.Range("BC2").Value = [double]"8.049493088558378e-06"
$ws.Range("BD2").Value = [double]"2.513934128945842e-08"
$ws.Range("BE2").Value = [double]"1.404703310082596e-13"
$ws.Range("BF2").Value = [double]"4.390376582856352e-14"
$ws.Range("BG2").Value = [double]"5.321237773533739e-16"
$ws.Range("BH2").Value = [double]"3.551311194456836e-17"
$ws.Range("BI2").Value = [double]"1.847624861396295e-17"
$ws.Range("BJ2").Value = [double]"1.09770537916695e-20"
$ws.Range("BK2").Value = [double]"1.947269496883026e-18"
$ws.Range("BL2").Value = [double]"1.982855054020422e-37"
$ws.Range("BM2").Value = [double]"6.912213102611216e-21"
$ws.Range("BN2").Value = [double]"9.234434769563976e-13"
$ws.Range("BO2").Value = [double]"1.111659093927002e-30"
$ws.Range("BP2").Value = [double]"4.390331159880988e-34"
$ws.Range("BQ2").Value = [double]"7.318558684780213e-11"
$ws.Range("BR2").Value = [double]"2.639221716934799e-36"
$ws.Range("BS2").Value = [double]"4.719125456714437e-11"
$ws.Range("BT2").Value = [double]"5.31034555661964e-26"
$ws.Range("BU2").Value = [double]"1.088924208331571e-12"
$ws.Range("BV2").Value = [double]"5.065814715513746e-14"
$ws.Range("BW2").Value = [double]"8.610568324938391e-07"
$ws.Range("BX2").Value = [double]"2.367709752548556e-07"
$ws.Range("BY2").Value = [double]"3.526589804781994e-13"
$ws.Range("CA2").Value = [double]"1.78854594291072e-15"
$ws.Range("CB2").Value = [double]"8.908427013128241e-17"
$ws.Range("CC2").Value = [double]"1.140648623239513e-19"
$ws.Range("CG2").Value = [double]"3.220858304440154e-07"
$ws.Range("CH2").Value = [double]"8.098475669025974e-05"
$ws.Range("CI2").Value = [double]"0.0004167856777646631"
$ws.Range("CJ2").Value = [double]"0.0007784035518922982"
$ws.Range("CK2").Value = [double]"0.001671908877086113"
$ws.Range("CL2").Value = [double]"0.006417731405726039"
$ws.Range("CM2").Value = [double]"99.90825478696034"
$ws.Range("CN2").Value = [double]"0.02092477296422002"
$ws.Range("CO2").Value = [double]"0.05877004516592325"
$ws.Range("CQ2").Value = [double]"0.001805729417158296"
$ws.Range("CR2").Value = [double]"0.0008771664155650092"
$ws.Range("CS2").Value = [double]"1.362612505605066e-06"
$ws.Range("CT2").Value = [double]"7.739728835086348e-11"
$ws.Range("CU2").Value = [double]"2.876045335280506e-11"
$ws.Range("CV2").Value = [double]"2.634737704517469e-12"
$ws.Range("CW2").Value = [double]"2.576365497228361e-13"
$ws.Range("CX2").Value = [double]"1.709191377736486e-13"
$ws.Range("CY2").Value = [double]"2.45932671091603e-16"
$ws.Range("CZ2").Value = [double]"1.123414406512061e-14"
$ws.Range("DA2").Value = [double]"2.145685625786699e-24"
$ws.Range("DB2").Value = [double]"1.320487653497826e-16"
$ws.Range("DC2").Value = [double]"1.739079842841245e-27"
$ws.Range("DD2").Value = [double]"8.69273901999894e-24"
$ws.Range("DE2").Value = [double]"1.745849505308884e-24"
$ws.Range("DF2").Value = [double]"4.803319652637677e-28"
$ws.Range("DG2").Value = [double]"7.287859653135067e-25"
$ws.Range("DH2").Value = [double]"7.794259718458143e-35"
$ws.Range("DI2").Value = [double]"6.000276146216289e-21"
$ws.Range("DJ2").Value = [double]"6.5806848350016e-14"
$ws.Range("DK2").Value = [double]"2.647201952027319e-15"
$ws.Range("DM2").Value = [double]"2.476483352287032e-31"
$ws.Range("DN2").Value = [double]"1.008004658302709e-14"
$ws.Range("DP2").Value = [double]"6.810955114454833e-17"
$ws.Range("DQ2").Value = [double]"4.062887552250365e-18"
$ws.Range("DR2").Value = [double]"2.863751518705569e-21"
$ws.Range("DV2").Value = [double]"3.199884010562363e-09"
$ws.Range("DW2").Value = [double]"1.722526977400943e-07"
$ws.Range("DX2").Value = [double]"2.771758454310595e-07"
$ws.Range("DY2").Value = [double]"8.045635810766119e-28"
$ws.Range("DZ2").Value = [double]"3.947241710938177e-20"
$ws.Range("EA2").Value = [double]"9.721426793360974e-13"
$ws.Range("EB2").Value = [double]"0.004616282061322658"
$ws.Range("EC2").Value = [double]"3.668595426600691e-07"
$ws.Range("ED2").Value = [double]"0.0001644347700471588"
$ws.Range("EF2").Value = [double]"0.01908600744327159"
$ws.Range("EG2").Value = [double]"94.92159406102468"
$ws.Range("EH2").Value = [double]"0.01342455317862403"
$ws.Range("EI2").Value = [double]"0.01090964834527415"
$ws.Range("EJ2").Value = [double]"0.008971747251799927"
$ws.Range("EK2").Value = [double]"2.241065384435269"
$ws.Range("EL2").Value = [double]"1.028518184711413"
$ws.Range("EM2").Value = [double]"1.582972917361153"
$ws.Range("EN2").Value = [double]"0.09428678177240127"
$ws.Range("EO2").Value = [double]"0.01961865402754561"
$ws.Range("EP2").Value = [double]"3.360348124239175e-08"
$ws.Range("EQ2").Value = [double]"0.04342824824604311"
$ws.Range("ER2").Value = [double]"5.831121389380853e-09"
$ws.Range("ES2").Value = [double]"1.662556139066817e-06"
$ws.Range("ET2").Value = [double]"6.219346331940888e-07"
$ws.Range("EU2").Value = [double]"4.96628978544766e-07"
$ws.Range("EV2").Value = [double]"2.741328608699776e-07"
$ws.Range("EW2").Value = [double]"3.043608611026193e-07"
$ws.Range("EX2").Value = [double]"0.003515079658799791"
$ws.Range("EY2").Value = [double]"6.637058635474689e-09"
$ws.Range("EZ2").Value = [double]"2.865680674594442e-10"
$ws.Range("FA2").Value = [double]"0.006236350659633702"
$ws.Range("FB2").Value = [double]"0.001587438543529343"
$ws.Range("FC2").Value = [double]"1.036959966395266e-09"
$ws.Range("FE2").Value = [double]"1.106631399766121e-11"
$ws.Range("FF2").Value = [double]"5.325176650829385e-13"
$ws.Range("FG2").Value = [double]"5.19857828500024e-16"
